# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
# Several match rows had their odds/result data mismatched to the wrong
# fixture (the unique match "id" in column B pointed at the wrong row's
# stats). This swaps the full data (columns B through AC, i.e. everything
# except the running row index in column A) between the affected row
# pairs so each row carries the correct id/HomeTeam/AwayTeam/odds set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2, $colStart, $colEnd) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# Column B = 2 ... Column AC = 29 (A = row index, left untouched)
$colStart = 2
$colEnd = 29

Swap-Rows 378 379 $colStart $colEnd
Swap-Rows 394 396 $colStart $colEnd
Swap-Rows 397 398 $colStart $colEnd
